$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.263.09"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.812.20"
$ws.Range("E3").Value = "  +3.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.73"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4324"
$ws.Range("E7").Value = "  -4.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.72"
$ws.Range("E9").Value = "  -1.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07630"
$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("E11").Value = "  +1.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.86"
$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.295"
$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.441"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.822.28"
$ws.Range("E16").Value = "  +3.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.27"
$ws.Range("E17").Value = "  +7.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001078"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06407"
$ws.Range("E19").Value = "  +2.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.37"
$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.208"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.261.91"
$ws.Range("E23").Value = "  +1.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.55"
$ws.Range("E24").Value = "  -0.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.137"
$ws.Range("E25").Value = "  -7.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.71"
$ws.Range("E26").Value = "  +4.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.60"
$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.027.69"
$ws.Range("E28").Value = "  +3.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.250"
$ws.Range("E29").Value = "  -4.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.29"
$ws.Range("E30").Value = "  +1.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.183"
$ws.Range("E31").Value = "  -2.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.970"
$ws.Range("E32").Value = "  +4.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09103"
$ws.Range("E33").Value = "  -2.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.541"
$ws.Range("E34").Value = "  -3.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.89"
$ws.Range("E35").Value = "  +1.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02387"
$ws.Range("E36").Value = "  +2.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.185"
$ws.Range("E37").Value = "  +1.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2157"
$ws.Range("E38").Value = "  -0.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6552"
$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06172"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.202"
$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.032"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.428"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.76"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6064"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.725"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.99"
$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.008"
$ws.Range("E49").Value = "  +0.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.170"
$ws.Range("E50").Value = "  +3.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06997"
$ws.Range("E51").Value = "  +1.12%  "
